$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, whether to force text format
$updates = @(
    @{Cell='D2'; Value='96.636.72'; ForceText=$false}
    @{Cell='E2'; Value='  -0.45%  '; ForceText=$false}
    @{Cell='D3'; Value='3.677.21'; ForceText=$false}
    @{Cell='E3'; Value='  +2.20%  '; ForceText=$false}
    @{Cell='E4'; Value='  -0.03%  '; ForceText=$false}
    @{Cell='D5'; Value='241.23'; ForceText=$true}
    @{Cell='E5'; Value='  -0.82%  '; ForceText=$false}
    @{Cell='E6'; Value='  +9.22%  '; ForceText=$false}
    @{Cell='D7'; Value='664.18'; ForceText=$true}
    @{Cell='E7'; Value='  +0.97%  '; ForceText=$false}
    @{Cell='D8'; Value='0.423'; ForceText=$true}
    @{Cell='E8'; Value='  +1.97%  '; ForceText=$false}
    @{Cell='E9'; Value='  +1.65%  '; ForceText=$false}
    @{Cell='E10'; Value='  +0.03%  '; ForceText=$false}
    @{Cell='D11'; Value='3.674.52'; ForceText=$false}
    @{Cell='E11'; Value='  +2.22%  '; ForceText=$false}
    @{Cell='D12'; Value='45.47'; ForceText=$true}
    @{Cell='E12'; Value='  +4.10%  '; ForceText=$false}
    @{Cell='E13'; Value='  +0.55%  '; ForceText=$false}
    @{Cell='D14'; Value='6.94'; ForceText=$true}
    @{Cell='E14'; Value='  +7.24%  '; ForceText=$false}
    @{Cell='D15'; Value='4.360.81'; ForceText=$false}
    @{Cell='E15'; Value='  +2.18%  '; ForceText=$false}
    @{Cell='E16'; Value='  +4.19%  '; ForceText=$false}
    @{Cell='D17'; Value='96.374.33'; ForceText=$false}
    @{Cell='E17'; Value='  -0.40%  '; ForceText=$false}
    @{Cell='D18'; Value='8.89'; ForceText=$true}
    @{Cell='E18'; Value='  +10.51%  '; ForceText=$false}
    @{Cell='D19'; Value='3.683.08'; ForceText=$false}
    @{Cell='E19'; Value='  +2.92%  '; ForceText=$false}
    @{Cell='D20'; Value='12.89'; ForceText=$true}
    @{Cell='E20'; Value='  +1.66%  '; ForceText=$false}
    @{Cell='D21'; Value='18.45'; ForceText=$true}
    @{Cell='E21'; Value='  +2.36%  '; ForceText=$false}
    @{Cell='D22'; Value='0.525'; ForceText=$true}
    @{Cell='E22'; Value='  -2.59%  '; ForceText=$false}
    @{Cell='D23'; Value='527.22'; ForceText=$true}
    @{Cell='E23'; Value='  +2.93%  '; ForceText=$false}
    @{Cell='D24'; Value='3.44'; ForceText=$true}
    @{Cell='E24'; Value='  +0.43%  '; ForceText=$false}
    @{Cell='E25'; Value='  +0.42%  '; ForceText=$false}
    @{Cell='D26'; Value='6.97'; ForceText=$true}
    @{Cell='E26'; Value='  +1.43%  '; ForceText=$false}
    @{Cell='D27'; Value='102.16'; ForceText=$true}
    @{Cell='D28'; Value='13.06'; ForceText=$true}
    @{Cell='E28'; Value='  -0.31%  '; ForceText=$false}
    @{Cell='B29'; Value='Hedera'; ForceText=$false}
    @{Cell='C29'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText=$false}
    @{Cell='D29'; Value='0.169'; ForceText=$true}
    @{Cell='E29'; Value='  +12.96%  '; ForceText=$false}
    @{Cell='B30'; Value='InternetComputer(DFINITY)'; ForceText=$false}
    @{Cell='C30'; Value='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText=$false}
    @{Cell='D30'; Value='12.57'; ForceText=$true}
    @{Cell='E30'; Value='  +8.26%  '; ForceText=$false}
    @{Cell='B31'; Value='PancakeSwap'; ForceText=$false}
    @{Cell='C31'; Value='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; ForceText=$false}
    @{Cell='D31'; Value='3.07'; ForceText=$true}
    @{Cell='E31'; Value='  +1.02%  '; ForceText=$false}
    @{Cell='B32'; Value='Dai'; ForceText=$false}
    @{Cell='C32'; Value='https://coinranking.com/coin/MoTuySvg7+dai-dai'; ForceText=$false}
    @{Cell='D32'; Value='1.00'; ForceText=$true}
    @{Cell='E32'; Value='  +0.09%  '; ForceText=$false}
    @{Cell='B33'; Value='Fetch.AI'; ForceText=$false}
    @{Cell='C33'; Value='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; ForceText=$false}
    @{Cell='D33'; Value='1.87'; ForceText=$true}
    @{Cell='E33'; Value='  +15.01%  '; ForceText=$false}
    @{Cell='D34'; Value='0.185'; ForceText=$true}
    @{Cell='E34'; Value='  -0.48%  '; ForceText=$false}
    @{Cell='B35'; Value='EthereumClassic'; ForceText=$false}
    @{Cell='C35'; Value='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; ForceText=$false}
    @{Cell='D35'; Value='32.82'; ForceText=$true}
    @{Cell='E35'; Value='  +3.57%  '; ForceText=$false}
    @{Cell='B36'; Value='Binance-PegBSC-USD'; ForceText=$false}
    @{Cell='C36'; Value='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; ForceText=$false}
    @{Cell='D36'; Value='1.00'; ForceText=$true}
    @{Cell='E36'; Value='  +0.20%  '; ForceText=$false}
    @{Cell='B37'; Value='PolygonEcosystemToken'; ForceText=$false}
    @{Cell='C37'; Value='https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'; ForceText=$false}
    @{Cell='D37'; Value='0.594'; ForceText=$true}
    @{Cell='E37'; Value='  +3.90%  '; ForceText=$false}
    @{Cell='B38'; Value='Bittensor'; ForceText=$false}
    @{Cell='C38'; Value='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; ForceText=$false}
    @{Cell='D38'; Value='636.53'; ForceText=$true}
    @{Cell='E38'; Value='  +2.21%  '; ForceText=$false}
    @{Cell='B39'; Value='RenderToken'; ForceText=$false}
    @{Cell='C39'; Value='https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'; ForceText=$false}
    @{Cell='D39'; Value='8.72'; ForceText=$true}
    @{Cell='E39'; Value='  -0.37%  '; ForceText=$false}
    @{Cell='B40'; Value='EnergySwap'; ForceText=$false}
    @{Cell='C40'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText=$false}
    @{Cell='D40'; Value='44.59'; ForceText=$true}
    @{Cell='E40'; Value='  +33.85%  '; ForceText=$false}
    @{Cell='B41'; Value='Kaspa'; ForceText=$false}
    @{Cell='C41'; Value='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; ForceText=$false}
    @{Cell='D41'; Value='0.161'; ForceText=$true}
    @{Cell='E41'; Value='  +5.43%  '; ForceText=$false}
    @{Cell='B42'; Value='ARBITRUM'; ForceText=$false}
    @{Cell='C42'; Value='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText=$false}
    @{Cell='D42'; Value='0.967'; ForceText=$true}
    @{Cell='E42'; Value='  +5.35%  '; ForceText=$false}
    @{Cell='B43'; Value='ImmutableX'; ForceText=$false}
    @{Cell='C43'; Value='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText=$false}
    @{Cell='D43'; Value='1.97'; ForceText=$true}
    @{Cell='E43'; Value='  +3.03%  '; ForceText=$false}
    @{Cell='B44'; Value='Filecoin'; ForceText=$false}
    @{Cell='C44'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText=$false}
    @{Cell='D44'; Value='6.43'; ForceText=$true}
    @{Cell='E44'; Value='  +8.61%  '; ForceText=$false}
    @{Cell='B45'; Value='USDe'; ForceText=$false}
    @{Cell='C45'; Value='https://coinranking.com/coin/exbfr2U-0+usde-usde'; ForceText=$false}
    @{Cell='D45'; Value='1.00'; ForceText=$true}
    @{Cell='E45'; Value='  +0.02%  '; ForceText=$false}
    @{Cell='B46'; Value='Algorand'; ForceText=$false}
    @{Cell='C46'; Value='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText=$false}
    @{Cell='D46'; Value='0.458'; ForceText=$true}
    @{Cell='E46'; Value='  +22.23%  '; ForceText=$false}
    @{Cell='B47'; Value='VeChain'; ForceText=$false}
    @{Cell='C47'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText=$false}
    @{Cell='D47'; Value='0.0461'; ForceText=$true}
    @{Cell='E47'; Value='  +6.74%  '; ForceText=$false}
    @{Cell='B48'; Value='Stacks'; ForceText=$false}
    @{Cell='C48'; Value='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; ForceText=$false}
    @{Cell='D48'; Value='2.29'; ForceText=$true}
    @{Cell='E48'; Value='  -0.92%  '; ForceText=$false}
    @{Cell='B49'; Value='WhiteBITCoin'; ForceText=$false}
    @{Cell='C49'; Value='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; ForceText=$false}
    @{Cell='D49'; Value='23.65'; ForceText=$true}
    @{Cell='E49'; Value='  -0.12%  '; ForceText=$false}
    @{Cell='B50'; Value='MantraDAO'; ForceText=$false}
    @{Cell='C50'; Value='https://coinranking.com/coin/cTdD8lD-6+mantradao-om'; ForceText=$false}
    @{Cell='D50'; Value='3.64'; ForceText=$true}
    @{Cell='E50'; Value='  +3.28%  '; ForceText=$false}
    @{Cell='D51'; Value='8.62'; ForceText=$true}
    @{Cell='E51'; Value='  +2.07%  '; ForceText=$false}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = '@'
    }
    $rng.Value = $u.Value
}
